# Implementacion de errores (Primera version)
#
# The lexical-analyzer transition matrix ("Matriz_Lexico") used shared string
# "q999" as a catch-all error/trap state. This change introduces a second,
# distinct error state "q998" and re-points a specific block of transitions
# (rows 86-91, the q84..q89 states) from "q999" to the new "q998" state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Contiguous ranges (per row) whose current value is "q999" and must become "q998".
$ranges = @(
    "B86:E86",
    "G86:AB86",
    "B87:AB87",
    "AM87:BG87",
    "B88:E88",
    "G88:AB88",
    "B89:AB89",
    "AM89:AV89",
    "AY89:BG89",
    "B90:AB90",
    "AM90:BG90",
    "B91:AB91"
)

foreach ($addr in $ranges) {
    $rng = $ws.Range($addr)
    foreach ($cell in $rng.Cells) {
        if ($cell.Value2 -eq "q999") {
            $cell.Value2 = "q998"
        }
    }
}

Write-Host "Updated $($ranges.Count) ranges (q999 -> q998)."

# Reflect the author's updated scroll/zoom/selection state on the sheet.
$ws.Activate()
$excel.ActiveWindow.Zoom = 87
[void]$ws.Range("A131").Select()
